$d = $word.ActiveDocument

# Locate the paragraph that starts the "site footer" block which was
# appended by the Jekyll build ("Ver no Jupiter Salvar em pdf Salvar em
# docx" followed by the "(c) 2020 ..." copyright line). That block, along
# with the blank paragraph immediately preceding it, is being removed by
# this commit; the blank paragraph that follows the copyright line (the
# one right before the page-break paragraph) must remain.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Ver no Jupiter*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Delete highest index first so the lower indices stay valid:
    #   targetIndex + 1 -> the "© 2020 ..." copyright paragraph
    #   targetIndex     -> the "Ver no Jupiter ..." paragraph
    #   targetIndex - 1 -> the blank paragraph right before it
    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()
    $d.Paragraphs.Item($targetIndex).Range.Delete()
    $d.Paragraphs.Item($targetIndex - 1).Range.Delete()
}
